$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row labels (column A) replacing the old asset class names
$labels = @(
    "Emerging Mkts",
    "US Treasuries",
    "High Yield",
    "Int'l Bonds",
    "Bonds - Agg",
    "Russ 1K Gro",
    "GOLD",
    "Commodities",
    "Small Stocks",
    "Russ 1K Val",
    "Real Estate",
    "S&P 500"
)

# New values for column B (Opt Portfolio) and C (Opt Portfolio with View)
# Written in plain decimal (no exponent) since the interpreter's numeric
# literal grammar does not accept scientific notation.
$bValues = @(
    0,
    0.3677972395916179,
    0.02992458714920683,
    0.1234796460533761,
    0.2978233882891922,
    0.01090413753029659,
    0,
    0.04514392892298908,
    0.04884331813316044,
    0.04509121052071238,
    0,
    0.03099254380944858
)

$cValues = @(
    0.0000000000000000004881338050412024,
    0.4552604038968207,
    0.06809171850192124,
    0,
    0.3660709093206801,
    0.0000000000000000009056290203122605,
    0.00000000000000002624961402989111,
    0.0511784746799404,
    0.02383451163665967,
    0.03556398196397785,
    0.0000000000000000004562268221530674,
    0.000000000000000000228076842600966
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

# Column A keeps the bold, bordered, centered label style used by the
# original rows - make sure the newly added rows (9-13) pick it up too.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2:A13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Apply percent number format to the newly populated B column range
$ws.Range("B2:B13").NumberFormat = "0.0%"

# Size columns to fit the new (longer) labels / percent values, matching
# the bestFit widths Excel computed for the author's data.
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666
$ws.Columns.Item(3).ColumnWidth = 21.666666666666668

# Page setup matching the re-saved workbook (paper size 9 = A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leftover selection cursor position from the author's last interactive edit
$ws.Range("B23").Select() | Out-Null

$wb.Save()
